$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C.
# Existing column C (Total Debt/Equity) shifts to D, existing column D (Price/BV) shifts to E.
$ws.Range("C1").EntireColumn.Insert()

# Give the new header cell (C1) the same look as the other header cells (bold/border/centered).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Header text for the newly inserted column.
$ws.Range("C1").Value = "Book Value [ExclRevalReserve]/Share (Rs.)"

# New column C data: Book Value [ExclRevalReserve]/Share (Rs.)
$cValues = @{
    2  = "708.19"
    3  = "697.01"
    4  = "736.20"
    5  = "617.15"
    6  = "639.41"
    7  = "496.66"
    8  = "886.83"
    9  = "741.28"
    10 = "667.98"
    11 = "609.76"
    12 = "554.17"
    13 = "498.22"
    14 = "446.30"
    15 = "392.51"
    16 = "727.78"
    17 = "542.83"
    18 = "439.67"
    19 = "324.11"
    20 = "270.43"
    21 = "227.22"
}

# Column D data (Total Debt/Equity (X)) - updated values after the shift.
$dValues = @{
    2  = "0.45"
    3  = "0.41"
    4  = "0.41"
    5  = "0.65"
    6  = "0.39"
    7  = "0.31"
    8  = "0.35"
    9  = "0.38"
    10 = "0.41"
    11 = "0.43"
    12 = "0.30"
    13 = "0.36"
    14 = "0.43"
    15 = "0.49"
    16 = "0.65"
    17 = "0.46"
    18 = "0.45"
    19 = "0.48"
    20 = "0.50"
    21 = "0.66"
}

foreach ($row in $cValues.Keys) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = "'" + $cValues[$row]
    $cell.Style = "Normal"
}

foreach ($row in $dValues.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = "'" + $dValues[$row]
    $cell.Style = "Normal"
}
